$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.67081492403963
$ws.Range("C2").Value = 4.920529886062094
$ws.Range("D2").Value = 7.999223543718954
$ws.Range("E2").Value = 10.17169003887623
$ws.Range("F2").Value = 39.13785774348442
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 32.23090148857823
$ws.Range("K2").Value = 13.68377707550442
$ws.Range("L2").Value = 10.40803754457147
$ws.Range("M2").Value = 15.90298431543091
$ws.Range("N2").Value = 23.27016272231272
$ws.Range("B3").Value = 15.51084573806452
$ws.Range("C3").Value = 4.652768619461452
$ws.Range("D3").Value = 8.004217635475538
$ws.Range("E3").Value = 10.18728975260467
$ws.Range("F3").Value = 39.08171923743588
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 32.26242961528182
$ws.Range("K3").Value = 13.56592547323157
$ws.Range("L3").Value = 10.4177235225461
$ws.Range("M3").Value = 15.89026421742978
$ws.Range("N3").Value = 23.32136920468524
$ws.Range("B4").Value = 15.41588974984399
$ws.Range("C4").Value = 4.479307964505583
$ws.Range("D4").Value = 8.007336526810194
$ws.Range("E4").Value = 10.19768113381797
$ws.Range("F4").Value = 39.05552495011177
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 32.28709074116836
$ws.Range("K4").Value = 13.49652669539663
$ws.Range("L4").Value = 10.42513447020941
$ws.Range("M4").Value = 15.88529043399297
$ws.Range("N4").Value = 23.35470868525073
$ws.Range("B5").Value = 15.37805888224988
$ws.Range("C5").Value = 4.406370221192327
$ws.Range("D5").Value = 8.008620847805496
$ws.Range("E5").Value = 10.20212056106498
$ws.Range("F5").Value = 39.04693701677829
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 32.2984725440927
$ws.Range("K5").Value = 13.46902004518511
$ws.Range("L5").Value = 10.42852285962393
$ws.Range("M5").Value = 15.88397935991165
$ws.Range("N5").Value = 23.3687727117918
$ws.Range("B6").Value = 15.37183053686751
$ws.Range("C6").Value = 4.394123872828688
$ws.Range("D6").Value = 8.008834918770644
$ws.Range("E6").Value = 10.20287010893341
$ws.Range("F6").Value = 39.04563714665337
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 32.30044290319917
$ws.Range("K6").Value = 13.46450014566105
$ws.Range("L6").Value = 10.42910775522113
$ws.Range("M6").Value = 15.88380495520753
$ws.Range("N6").Value = 23.37113691341567
$ws.Range("B7").Value = 15.41537599349619
$ws.Range("C7").Value = 4.478333376939352
$ws.Range("D7").Value = 8.007353793388146
$ws.Range("E7").Value = 10.19774017550424
$ws.Range("F7").Value = 39.05540067610969
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 32.28723884804788
$ws.Range("K7").Value = 13.49615256043114
$ws.Range("L7").Value = 10.42517867535337
$ws.Range("M7").Value = 15.88526985115359
$ws.Range("N7").Value = 23.35489642189018
$ws.Range("B8").Value = 15.61500622883879
$ws.Range("C8").Value = 4.830098539612433
$ws.Range("D8").Value = 8.000934699060853
$ws.Range("E8").Value = 10.17690030058465
$ws.Range("F8").Value = 39.11678756644253
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 32.24067115670103
$ws.Range("K8").Value = 13.64254477777395
$ws.Range("L8").Value = 10.41107367539833
$ws.Range("M8").Value = 15.89801136932094
$ws.Range("N8").Value = 23.28742500113151
$ws.Range("B9").Value = 16.03041317888741
$ws.Range("C9").Value = 5.447225721338106
$ws.Range("D9").Value = 7.988756511359109
$ws.Range("E9").Value = 10.14246741738077
$ws.Range("F9").Value = 39.30251150559797
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 32.19147216718096
$ws.Range("K9").Value = 13.95178599548914
$ws.Range("L9").Value = 10.39501275510289
$ws.Range("M9").Value = 15.94537612566934
$ws.Range("N9").Value = 23.17015333071659
$ws.Range("B10").Value = 16.3474652671117
$ws.Range("C10").Value = 5.855514113362712
$ws.Range("D10").Value = 7.980048883276697
$ws.Range("E10").Value = 10.12106829534417
$ws.Range("F10").Value = 39.47828208499528
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 32.18105145850652
$ws.Range("K10").Value = 14.1906215500092
$ws.Range("L10").Value = 10.39026118279253
$ws.Range("M10").Value = 15.99362533723024
$ws.Range("N10").Value = 23.09312848792216
$ws.Range("B11").Value = 16.493676756344
$ws.Range("C11").Value = 6.031360694796543
$ws.Range("D11").Value = 7.976137439965781
$ws.Range("E11").Value = 10.11217483216902
$ws.Range("F11").Value = 39.56663885232372
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 32.18190047364136
$ws.Range("K11").Value = 14.30138320636045
$ws.Range("L11").Value = 10.38962382720502
$ws.Range("M11").Value = 16.01844655295969
$ws.Range("N11").Value = 23.06006447103063
$ws.Range("B12").Value = 16.54927801328267
$ws.Range("C12").Value = 6.096522621389966
$ws.Range("D12").Value = 7.974663265944256
$ws.Range("E12").Value = 10.10892765713576
$ws.Range("F12").Value = 39.60128997984649
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 32.18302528627884
$ws.Range("K12").Value = 14.34359398046475
$ws.Range("L12").Value = 10.38960096941503
$ws.Range("M12").Value = 16.02825388437261
$ws.Range("N12").Value = 23.04782753356908
$ws.Range("B13").Value = 16.53729368363017
$ws.Range("C13").Value = 6.08255238613914
$ws.Range("D13").Value = 7.974980446411682
$ws.Range("E13").Value = 10.10962163705703
$ws.Range("F13").Value = 39.59377446976396
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 32.18274732152715
$ws.Range("K13").Value = 14.33449180699852
$ws.Range("L13").Value = 10.38959618626666
$ws.Range("M13").Value = 16.02612363420488
$ws.Range("N13").Value = 23.05045036462083
$ws.Range("B14").Value = 16.49824670394947
$ws.Range("C14").Value = 6.036750237594812
$ws.Range("D14").Value = 7.976016019210021
$ws.Range("E14").Value = 10.11190527047494
$ws.Range("F14").Value = 39.56946581824784
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 32.18197691706253
$ws.Range("K14").Value = 14.30485075479252
$ws.Range("L14").Value = 10.38961757188096
$ws.Range("M14").Value = 16.01924525946925
$ws.Range("N14").Value = 23.0590520469995
$ws.Range("B15").Value = 16.47435828370187
$ws.Range("C15").Value = 6.008509080215854
$ws.Range("D15").Value = 7.976651245459429
$ws.Range("E15").Value = 10.11331975500774
$ws.Range("F15").Value = 39.55473086490601
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 32.18160961685638
$ws.Range("K15").Value = 14.28672854159029
$ws.Range("L15").Value = 10.38965910421002
$ws.Range("M15").Value = 16.01508504545805
$ws.Range("N15").Value = 23.06435776150711
$ws.Range("B16").Value = 16.33794553951174
$ws.Range("C16").Value = 5.843822652249844
$ws.Range("D16").Value = 7.980305493217377
$ws.Range("E16").Value = 10.12166639820321
$ws.Range("F16").Value = 39.47267548008035
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 32.18110841784939
$ws.Range("K16").Value = 14.18342246244556
$ws.Range("L16").Value = 10.39033344938797
$ws.Range("M16").Value = 15.99206059660211
$ws.Range("N16").Value = 23.09532901369095
$ws.Range("B17").Value = 16.25473227284342
$ws.Range("C17").Value = 5.740256028934073
$ws.Range("D17").Value = 7.982559883858029
$ws.Range("E17").Value = 10.12700196212954
$ws.Range("F17").Value = 39.42447769505466
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 32.18223234306664
$ws.Range("K17").Value = 14.12056329490604
$ws.Range("L17").Value = 10.39113705548441
$ws.Range("M17").Value = 15.97866812432373
$ws.Range("N17").Value = 23.11483448443402
$ws.Range("B18").Value = 16.20705989357798
$ws.Range("C18").Value = 5.679756932962062
$ws.Range("D18").Value = 7.983861236152974
$ws.Range("E18").Value = 10.13015002791398
$ws.Range("F18").Value = 39.39754706446431
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 32.18340495094271
$ws.Range("K18").Value = 14.08460955379606
$ws.Range("L18").Value = 10.39174277504116
$ws.Range("M18").Value = 15.97123581623391
$ws.Range("N18").Value = 23.12623942610166
$ws.Range("B19").Value = 16.19095298147172
$ws.Range("C19").Value = 5.659113409033678
$ws.Range("D19").Value = 7.98430266094023
$ws.Range("E19").Value = 10.13122952100779
$ws.Range("F19").Value = 39.38856518894121
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 32.18389235104092
$ws.Range("K19").Value = 14.07247189031437
$ws.Range("L19").Value = 10.39197252799209
$ws.Range("M19").Value = 15.96876600028644
$ws.Range("N19").Value = 23.13013288354616
$ws.Range("B20").Value = 16.26357120492075
$ws.Range("C20").Value = 5.751377205551946
$ws.Range("D20").Value = 7.982319416118644
$ws.Range("E20").Value = 10.12642578895514
$ws.Range("F20").Value = 39.42952662381807
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 32.18205824612645
$ws.Range("K20").Value = 14.12723420023634
$ws.Range("L20").Value = 10.39103666250139
$ws.Range("M20").Value = 15.98006579544212
$ws.Range("N20").Value = 23.11273885497573
$ws.Range("B21").Value = 16.5097097991256
$ws.Range("C21").Value = 6.05024221241837
$ws.Range("D21").Value = 7.975711657398916
$ws.Range("E21").Value = 10.11123124220103
$ws.Range("F21").Value = 39.57657362819641
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 32.18218140740156
$ws.Range("K21").Value = 14.31355007411592
$ws.Range("L21").Value = 10.38960536631324
$ws.Range("M21").Value = 16.02125457140702
$ws.Range("N21").Value = 23.05651782753882
$ws.Range("B22").Value = 16.67191692635529
$ws.Range("C22").Value = 6.237246793382304
$ws.Range("D22").Value = 7.971433868373263
$ws.Range("E22").Value = 10.10200342875458
$ws.Range("F22").Value = 39.679619020225
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 32.18694372270032
$ws.Range("K22").Value = 14.43686192659677
$ws.Range("L22").Value = 10.38994311624116
$ws.Range("M22").Value = 16.05054981002324
$ws.Range("N22").Value = 23.02142755720866
$ws.Range("B23").Value = 16.58523796788667
$ws.Range("C23").Value = 6.138201672313301
$ws.Range("D23").Value = 7.973713321792995
$ws.Range("E23").Value = 10.1068643094124
$ws.Range("F23").Value = 39.62399203833334
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 32.18397384482444
$ws.Range("K23").Value = 14.37091885319942
$ws.Range("L23").Value = 10.38964660076762
$ws.Range("M23").Value = 16.03469873142138
$ws.Range("N23").Value = 23.04000470438888
$ws.Range("B24").Value = 16.25957459640254
$ws.Range("C24").Value = 5.746352302898413
$ws.Range("D24").Value = 7.982428115122573
$ws.Range("E24").Value = 10.12668602583533
$ws.Range("F24").Value = 39.42724157559492
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 32.18213531537764
$ws.Range("K24").Value = 14.12421770542455
$ws.Range("L24").Value = 10.39108160245739
$ws.Range("M24").Value = 15.97943307563033
$ws.Range("N24").Value = 23.11368569405315
$ws.Range("B25").Value = 15.91575153478295
$ws.Range("C25").Value = 5.288155719959859
$ws.Range("D25").Value = 7.99200822099141
$ws.Range("E25").Value = 10.1510960593542
$ws.Range("F25").Value = 39.2453194720651
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 32.20026618319334
$ws.Range("K25").Value = 13.86594494381697
$ws.Range("L25").Value = 10.39811795454702
$ws.Range("M25").Value = 15.93018514185313
$ws.Range("N25").Value = 23.20027190510183
